# Update the "展览" (Exhibition) and "全部类型" (All types) sheets:
#   - F2 (想去人数 / "want to go" count) for the first event changes 1306 -> 1339
#   - A new row 4 is appended describing a third event
#
# Both sheets carry identical data in this workbook, so the same edits are
# applied to each of them.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Update the "want to go" count on the existing first event row.
    $ws.Range("F2").Value = 1339

    # Append a new row (row 4) for the third event. Copy the formatting from
    # row 3 first (keeps the bold/bordered/centered style on column A, same
    # as A2/A3), then fill in the new values.
    $ws.Range("A3:I3").Copy($ws.Range("A4:I4"))

    $ws.Range("A4").Value = 3

    # B4 holds a plain "2024-05-18" date string (same as B2/B3), not a real
    # date value - format the cell as Text first so Excel doesn't silently
    # convert the ISO-looking string into a date serial number.
    $ws.Range("B4").NumberFormat = "@"
    $ws.Range("B4").Value = "2024-05-18"

    $ws.Range("C4").Value = "丽水·第三届HP国风动漫游戏嘉年华"
    $ws.Range("D4").Value = "好溪路与望城路交汇西北侧地块 丽水市水上运动中心"
    $ws.Range("E4").Value = "2024.05.18 09:00-05.18 17:00"
    $ws.Range("F4").Value = 0
    $ws.Range("G4").Value = 68
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=82901"
    $ws.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202403/sl5TubQI1710410535537.jpeg"
}

Write-Output "done"
